$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update phone numbers in B2 and B3
$ws.Range("B2").Value = 5532999999999
$ws.Range("B3").Value = 5532999999999

# Update selection to E8
$ws.Range("E8").Select()
